$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.056.80'
$ws.Range('D3').Value = '2.050.39'
$ws.Range('E3').Value = '  +0.22%  '
$ws.Range('E4').Value = '  +0.14%  '
$ws.Range('D5').Value = "'247.01"
$ws.Range('E5').Value = '  -0.42%  '
$ws.Range('D6').Value = "'0.662"
$ws.Range('E6').Value = '  +0.02%  '
$ws.Range('B7').Value = 'Solana'
$ws.Range('C7').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
$ws.Range('D7').Value = "'57.56"
$ws.Range('E7').Value = '  +1.97%  '
$ws.Range('B8').Value = 'USDC'
$ws.Range('C8').Value = 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc'
$ws.Range('D8').Value = "'1.00"
$ws.Range('E8').Value = '  -0.07%  '
$ws.Range('E9').Value = '  +0.68%  '
$ws.Range('D10').Value = "'0.0776"
$ws.Range('E10').Value = '  -0.34%  '
$ws.Range('E11').Value = '  +2.05%  '
$ws.Range('D12').Value = "'15.74"
$ws.Range('E12').Value = '  -0.04%  '
$ws.Range('D13').Value = "'0.901"
$ws.Range('E13').Value = '  +14.18%  '
$ws.Range('D14').Value = '2.350.72'
$ws.Range('E14').Value = '  +0.32%  '
$ws.Range('E15').Value = '  +2.99%  '
$ws.Range('D16').Value = '2.053.08'
$ws.Range('E16').Value = '  +0.28%  '
$ws.Range('D17').Value = "'18.49"
$ws.Range('E17').Value = '  +14.23%  '
$ws.Range('D18').Value = '37.023.47'
$ws.Range('E18').Value = '  +0.30%  '
$ws.Range('D19').Value = "'74.62"
$ws.Range('E19').Value = '  +0.81%  '
$ws.Range('D20').Value = '0.0₃0900'
$ws.Range('E20').Value = '  +0.84%  '
$ws.Range('D21').Value = "'5.49"
$ws.Range('E21').Value = '  +3.30%  '
$ws.Range('D22').Value = "'236.61"
$ws.Range('E22').Value = '  +0.30%  '
$ws.Range('E23').Value = '  -0.02%  '
$ws.Range('D24').Value = "'2.47"
$ws.Range('E24').Value = '  +5.04%  '
$ws.Range('D25').Value = "'9.54"
$ws.Range('E25').Value = '  +5.49%  '
$ws.Range('D26').Value = "'170.39"
$ws.Range('E26').Value = '  +1.62%  '
$ws.Range('D27').Value = "'2.17"
$ws.Range('E27').Value = '  -1.60%  '
$ws.Range('D28').Value = "'20.06"
$ws.Range('E28').Value = '  +1.39%  '
$ws.Range('E29').Value = '  +15.95%  '
$ws.Range('E30').Value = '  +0.04%  '
$ws.Range('E31').Value = '  +2.77%  '
$ws.Range('E32').Value = '  +9.29%  '
$ws.Range('D33').Value = "'0.0621"
$ws.Range('E33').Value = '  +1.47%  '
$ws.Range('D34').Value = "'0.0877"
$ws.Range('E34').Value = '  +0.44%  '
$ws.Range('E35').Value = '  +0.01%  '
$ws.Range('D36').Value = "'2.30"
$ws.Range('E36').Value = '  +3.82%  '
$ws.Range('E37').Value = '  +5.14%  '
$ws.Range('D38').Value = "'1.33"
$ws.Range('E38').Value = '  -0.61%  '
$ws.Range('D39').Value = "'3.10"
$ws.Range('E39').Value = '  -2.32%  '
$ws.Range('D40').Value = "'5.14"
$ws.Range('E40').Value = '  +4.50%  '
$ws.Range('D41').Value = "'0.0999"
$ws.Range('E41').Value = '  -5.35%  '
$ws.Range('E42').Value = '  +1.45%  '
$ws.Range('D43').Value = "'1.16"
$ws.Range('E43').Value = '  +3.90%  '
$ws.Range('D44').Value = "'98.83"
$ws.Range('E44').Value = '  +3.61%  '
$ws.Range('D45').Value = "'17.13"
$ws.Range('E45').Value = '  -0.52%  '
$ws.Range('D46').Value = "'2.39"
$ws.Range('E46').Value = '  -1.44%  '
$ws.Range('D47').Value = '1.301.10'
$ws.Range('E47').Value = '  +1.79%  '
$ws.Range('E48').Value = '  +0.92%  '
$ws.Range('D49').Value = "'6.84"
$ws.Range('E49').Value = '  +2.71%  '
$ws.Range('E50').Value = '  +5.01%  '
$ws.Range('D51').Value = '2.236.49'
$ws.Range('E51').Value = '  +0.19%  '
